$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'28.016.15"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.41%  '

$c = $ws.Range("D3")
$c.Value = "'1.868.70"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.82%  '

$c = $ws.Range("D4")
$c.Value = "'1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '

$c = $ws.Range("D5")
$c.Value = "'312.28"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '

$c = $ws.Range("D6")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '

$c = $ws.Range("D7")
$c.Value = "'0.5113"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.67%  '

$c = $ws.Range("D8")
$c.Value = "'0.3872"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.26%  '

$c = $ws.Range("D9")
$c.Value = "'0.08338"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '

$c = $ws.Range("D10")
$c.Value = "'1.111"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '

$c = $ws.Range("D11")
$c.Value = "'41.50"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.35%  '

$c = $ws.Range("D12")
$c.Value = "'6.163"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.62%  '

$c = $ws.Range("D13")
$c.Value = "'1.873.83"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.13%  '

$c = $ws.Range("D14")
$c.Value = "'20.47"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '

$c = $ws.Range("D15")
$c.Value = "'7.258"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.64%  '

$c = $ws.Range("D16")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '

$c = $ws.Range("D17")
$c.Value = "'0.00001098"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.18%  '

$c = $ws.Range("D18")
$c.Value = "'90.74"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.47%  '

$c = $ws.Range("D19")
$c.Value = "'0.06623"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.33%  '

$c = $ws.Range("D20")
$c.Value = "'17.60"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.76%  '

$c = $ws.Range("D21")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.23%  '

$c = $ws.Range("D22")
$c.Value = "'6.006"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.40%  '

$c = $ws.Range("D23")
$c.Value = "'28.073.50"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.34%  '

$c = $ws.Range("D24")
$c.Value = "'11.08"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '

$c = $ws.Range("D25")
$c.Value = "'2.244"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.17%  '

$c = $ws.Range("D26")
$c.Value = "'2.080.83"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.58%  '

$c = $ws.Range("D27")
$c.Value = "'2.464"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.34%  '

$c = $ws.Range("D28")
$c.Value = "'158.05"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.25%  '

$c = $ws.Range("D29")
$c.Value = "'20.49"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.00%  '

$c = $ws.Range("D30")
$c.Value = "'124.60"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.31%  '

$c = $ws.Range("D31")
$c.Value = "'0.1061"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.75%  '

$c = $ws.Range("D32")
$c.Value = "'1.031"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.48%  '

$c = $ws.Range("D33")
$c.Value = "'5.857"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.97%  '

$c = $ws.Range("D34")
$c.Value = "'3.599"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.30%  '

$c = $ws.Range("D35")
$c.Value = "'9.395"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.11%  '

$c = $ws.Range("D36")
$c.Value = "'0.02429"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.05%  '

$c = $ws.Range("D37")
$c.Value = "'0.06528"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.09%  '

$c = $ws.Range("D38")
$c.Value = "'0.2179"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '

$c = $ws.Range("D39")
$c.Value = "'1.200"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.27%  '

$c = $ws.Range("D40")
$c.Value = "'0.6469"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.59%  '

$c = $ws.Range("D41")
$c.Value = "'4.992"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.19%  '

$c = $ws.Range("D42")
$c.Value = "'1.218"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.34%  '

$c = $ws.Range("D43")
$c.Value = "'11.29"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.81%  '

$c = $ws.Range("D44")
$c.Value = "'0.6067"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.66%  '

$c = $ws.Range("D45")
$c.Value = "'12.91"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.77%  '

$c = $ws.Range("D46")
$c.Value = "'1.280"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '

$c = $ws.Range("D47")
$c.Value = "'3.673"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.19%  '

$c = $ws.Range("D48")
$c.Value = "'2.003"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.37%  '

$c = $ws.Range("D49")
$c.Value = "'1.214"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.37%  '

$c = $ws.Range("D50")
$c.Value = "'120.83"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

$c = $ws.Range("D51")
$c.Value = "'77.88"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.28%  '
